$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows where only values changed ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.767.44"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.91"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.86"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.42"
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.367"
$ws.Range("E7").Value = "  -10.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.891"
$ws.Range("E8").Value = "  +23.14%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.086.46"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.665"
$ws.Range("E11").Value = "  +17.46%  "
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -7.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.36"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.964.18"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.52"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.674.82"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.106.50"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000210"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.54"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.23"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.02"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.48"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.02"
$ws.Range("E26").Value = "  +6.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.265.64"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +9.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.166"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.16"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "508.06"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.61"
$ws.Range("E34").Value = "  -11.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.96"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.30"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.142"
$ws.Range("E43").Value = "  +13.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.372"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.84"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.22"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "157.06"
$ws.Range("E51").Value = "  -7.21%  "

# --- Rows 47-50 were reordered (coins swapped) with updated price/volume ---
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.80"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0687"
$ws.Range("E48").Value = "  +11.49%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.09"
$ws.Range("E49").Value = "  +2.03%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.22"
$ws.Range("E50").Value = "  +1.26%  "